# Generate Report for Handoff
# Updates the localization-status report so that the previously
# "Handed back" b.md row now reflects that it is "Ready for handoff"
# again (its handback content was found stale / not the latest),
# on the Overview sheet as well as on each per-locale sheet
# (zh-cn, de-de).

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item(1)   # Overview
$ws2 = $wb.Worksheets.Item(2)   # zh-cn
$ws3 = $wb.Worksheets.Item(3)   # de-de

$errorDetail = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/66114aee9abbaade7769563da60cbea8375cbbcc/e2e/a.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/e5003043ee47941c7b35f14bf5de1c0415f07dc6/e2e/b.md."

# ---- Overview sheet: row for b.md (row 3) ----
$ws1.Range("E3").Value2 = "Ready for handoff"
$ws1.Range("F3").Value2 = "Ready for handoff"
$ws1.Range("G3").Value2 = "2016-08-24 20:38:47"

# ---- zh-cn sheet: row for b.md (row 3) ----
$ws2.Range("C3").Value2 = "Ready for handoff"
$ws2.Range("F3").Value2 = "'False"
$ws2.Range("F3").Style = "Normal"
$ws2.Range("G3").Value2 = "b.63290e5768f688058c7b37413b0a5c26c308f864.zh-cn.xlf"
$ws2.Range("H3").Value2 = "2016-08-24 20:38:43"
$ws2.Range("P3").Value2 = $errorDetail
$ws2.Range("P1").ColumnWidth = 39.17

# ---- de-de sheet: row for b.md (row 3) ----
$ws3.Range("C3").Value2 = "Ready for handoff"
$ws3.Range("F3").Value2 = "'False"
$ws3.Range("F3").Style = "Normal"
$ws3.Range("G3").Value2 = "b.63290e5768f688058c7b37413b0a5c26c308f864.de-de.xlf"
$ws3.Range("H3").Value2 = "2016-08-24 20:38:47"
$ws3.Range("P3").Value2 = $errorDetail
$ws3.Range("P1").ColumnWidth = 39.17
